# Update Team Member Roles
#
# Slide 1, Shape "TextBox 4" holds the bulleted role list. Each bullet is a
# paragraph made of two runs ("<role> " + "Manager/Developer – <names>").
# This edit:
#   1. Re-splits paragraph 1's runs ("Team Leader/Project Manager – " / "Eliseo ")
#   2. Re-types paragraphs 2-6 in place (same text) so their two runs coalesce
#      into a single run
#   3. Re-types paragraph 7 in place (same text) so its two runs coalesce
#      into a single run (its end-of-paragraph mark is left untouched)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$tb = $s.Shapes.Item(2)

# --- paragraphs 2-7: collapse the two runs into one by re-typing the
#     paragraph's own (unchanged) text across the old run boundary -------
$paragraphSpans = @(
    @{ Start = 38;  Length = 40; Text = "Client-side Developer – Jennifer, Andrea" },
    @{ Start = 79;  Length = 46; Text = "Server-side Developer – Jennifer, Alex, Eliseo" },
    @{ Start = 126; Length = 32; Text = "Database Manager – Eliseo, Clark" },
    @{ Start = 159; Length = 37; Text = "UI/UX Design Manager – Alex, Jennifer" },
    @{ Start = 197; Length = 28; Text = "Test Manager – Eliseo, Clark" },
    @{ Start = 226; Length = 53; Text = "Business Analyst/Requirements Manager – Andrea, Clark" }
)

foreach ($span in $paragraphSpans) {
    $run = $tb.TextFrame.TextRange.Characters($span.Start, $span.Length)
    $run.Text = $span.Text
}

# --- paragraph 1: re-split "Team Leader/Project " / "Manager – Eliseo"
#     into "Team Leader/Project Manager – " / "Eliseo " -------------------
$firstPart = $tb.TextFrame.TextRange.Characters(1, 30)
$firstPart.Text = "Team Leader/Project Manager – "

$secondPart = $tb.TextFrame.TextRange.Characters(31, 6)
$secondPart.Text = "Eliseo "
